$d = $word.ActiveDocument

$pairs = @(
    @("154×6=924", "205×9=1845"),
    @("483×9=4347", "594×9=5346"),
    @("537×7=3759", "607×6=3642"),
    @("660×6=3960", "602×3=1806"),
    @("914×4=3656", "346×2=692"),
    @("843×2=1686", "743×4=2972"),
    @("661×6=3966", "565×3=1695"),
    @("837×4=3348", "871×5=4355"),
    @("230×7=1610", "597×6=3582"),
    @("123×9=1107", "812×2=1624"),
    @("856×6=5136", "695×4=2780"),
    @("137×3=411", "103×5=515"),
    @("458×9=4122", "545×8=4360"),
    @("615×6=3690", "575×4=2300"),
    @("510×4=2040", "349×3=1047"),
    @("660×5=3300", "692×2=1384"),
    @("642×9=5778", "378×3=1134"),
    @("857×5=4285", "590×3=1770"),
    @("117×8=936", "766×9=6894"),
    @("447×2=894", "169×6=1014"),
    @("665×2=1330", "746×8=5968"),
    @("582×8=4656", "995×3=2985"),
    @("910×6=5460", "403×4=1612"),
    @("416×4=1664", "152×7=1064"),
    @("663×3=1989", "911×5=4555")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
